# Realestate Update resale numbers 2024-01-19 17:55
# Appends a new data row (row 76) to the CityResaleNum sheet with the
# latest resale-number snapshot, matching the style of the existing rows
# (plain text for Date/Time/Weekday/Week, plain numbers for the city
# columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 76

# --- Text columns (A:D) -----------------------------------------------
# Values such as "2024-01-19", "17:55:57" or "02" would otherwise be
# auto-converted by Excel into a date/time serial or a plain number,
# which would both change their value and drop the leading zero. Format
# the cells as Text first so the literal string is preserved, then clear
# the formatting again afterwards so the cells end up unstyled, exactly
# like the rest of the data rows.
$textRange = $ws.Range("A$row`:D$row")
$textRange.NumberFormat = "@"

$ws.Range("A$row").Value = "2024-01-19"
$ws.Range("B$row").Value = "17:55:57"
$ws.Range("C$row").Value = "Friday"
$ws.Range("D$row").Value = "02"

$textRange.ClearFormats()

# --- Numeric columns (E:T) ---------------------------------------------
$ws.Range("E$row").Value = 138010
$ws.Range("F$row").Value = 140483
$ws.Range("G$row").Value = 171506
$ws.Range("H$row").Value = 148837
$ws.Range("I$row").Value = -1
$ws.Range("J$row").Value = 122319
$ws.Range("K$row").Value = 223585
$ws.Range("L$row").Value = 255136
$ws.Range("M$row").Value = 185268
$ws.Range("N$row").Value = 110317
$ws.Range("O$row").Value = 41384
$ws.Range("P$row").Value = 30906
$ws.Range("Q$row").Value = 73595
$ws.Range("R$row").Value = -1
$ws.Range("S$row").Value = 42717
$ws.Range("T$row").Value = -1
